$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 1252
$ws.Range("J55").Value = 1584.1666
$ws.Range("L55").Value = 1584.1666
$ws.Range("N55").Value = -2012.1666

$ws.Range("H74").Value = 6071.3125
$ws.Range("I74").Value = 3357
$ws.Range("K74").Value = 3357
$ws.Range("M74").Value = -2421

$ws.Range("H77").Value = 6071.3125
$ws.Range("I77").Value = 3357
$ws.Range("K77").Value = 16785
$ws.Range("M77").Value = -12105

$ws.Range("H88").Value = 11900.6
$ws.Range("I88").Value = 10499.667
$ws.Range("J88").Value = 14002
$ws.Range("K88").Value = 10499.667
$ws.Range("L88").Value = 14002
$ws.Range("M88").Value = -10093.667
$ws.Range("N88").Value = -14814

$ws.Range("H91").Value = 11900.6
$ws.Range("I91").Value = 10499.667
$ws.Range("J91").Value = 14002
$ws.Range("K91").Value = 10499.667
$ws.Range("L91").Value = 14002
$ws.Range("M91").Value = -9095.666999999999
$ws.Range("N91").Value = -16810

$ws.Range("H96").Value = 777.2
$ws.Range("I96").Value = 777.2
$ws.Range("K96").Value = 2331.6
$ws.Range("M96").Value = -958.6000000000004

$ws.Range("H99").Value = 184
$ws.Range("I99").Value = 184
$ws.Range("K99").Value = 552
$ws.Range("M99").Value = 946

$ws.Range("H103").Value = 433.57574
$ws.Range("J103").Value = 557.4
$ws.Range("L103").Value = 1672.2
$ws.Range("N103").Value = -2844.2

$ws.Range("H127").Value = 14528.866
$ws.Range("I127").Value = 15693.538
$ws.Range("K127").Value = 47080.614
$ws.Range("M127").Value = -42120.614

$ws.Range("H132").Value = 1963.2972
$ws.Range("I132").Value = 1284.7333
$ws.Range("J132").Value = 4871.4287
$ws.Range("K132").Value = 3854.199900000001
$ws.Range("L132").Value = 14614.2861
$ws.Range("M132").Value = -1324.199900000001
$ws.Range("N132").Value = -19674.2861

$ws.Range("H138").Value = 5670.074
$ws.Range("I138").Value = 3092.5
$ws.Range("J138").Value = 6958.8613
$ws.Range("K138").Value = 9277.5
$ws.Range("L138").Value = 20876.5839
$ws.Range("M138").Value = -4137.5
$ws.Range("N138").Value = -31156.5839

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 8222.929
$ws.Range("I45").Value = 4589
$ws.Range("K45").Value = 4589
$ws.Range("M45").Value = -4212

$ws.Range("H57").Value = 9599.799999999999
$ws.Range("I57").Value = 9599.799999999999
$ws.Range("K57").Value = 9599.799999999999
$ws.Range("M57").Value = -9115.799999999999

$ws.Range("H88").Value = 2939.2222
$ws.Range("I88").Value = 4045
$ws.Range("J88").Value = 2054.6
$ws.Range("K88").Value = 4045
$ws.Range("L88").Value = 2054.6
$ws.Range("M88").Value = -3639
$ws.Range("N88").Value = -2866.6

$ws.Range("H91").Value = 2939.2222
$ws.Range("I91").Value = 4045
$ws.Range("J91").Value = 2054.6
$ws.Range("K91").Value = 4045
$ws.Range("L91").Value = 2054.6
$ws.Range("M91").Value = -2641
$ws.Range("N91").Value = -4862.6

$ws.Range("H97").Value = 2186.5
$ws.Range("I97").Value = 2437.2
$ws.Range("J97").Value = 933
$ws.Range("K97").Value = 2437.2
$ws.Range("L97").Value = 933
$ws.Range("M97").Value = -1941.2
$ws.Range("N97").Value = -1925

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 851590.9399999999
$ws.Range("I86").Value = 1063809.4
$ws.Range("J86").Value = 2717.25
$ws.Range("K86").Value = 1063809.4
$ws.Range("L86").Value = 2717.25
$ws.Range("M86").Value = -1062686.4
$ws.Range("N86").Value = -4963.25

$ws.Range("H89").Value = 851590.9399999999
$ws.Range("I89").Value = 1063809.4
$ws.Range("J89").Value = 2717.25
$ws.Range("K89").Value = 5319047
$ws.Range("L89").Value = 13586.25
$ws.Range("M89").Value = -5313431
$ws.Range("N89").Value = -24818.25

$ws.Range("H94").Value = 1109.2
$ws.Range("I94").Value = 1124
$ws.Range("J94").Value = 1050
$ws.Range("K94").Value = 1124
$ws.Range("L94").Value = 1050
$ws.Range("M94").Value = -673
$ws.Range("N94").Value = -1952

$ws.Range("H109").Value = 37500
$ws.Range("J109").Value = 37500
$ws.Range("L109").Value = 37500
$ws.Range("N109").Value = -40274

$ws.Range("H126").Value = 84325.664
$ws.Range("J126").Value = 84325.664
$ws.Range("L126").Value = 84325.664
$ws.Range("N126").Value = -94205.664

$ws.Range("H128").Value = 12500
$ws.Range("I128").Value = 12500
$ws.Range("K128").Value = 37500
$ws.Range("M128").Value = -35010

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1210
$ws.Range("I22").Value = 294
$ws.Range("K22").Value = 294
$ws.Range("M22").Value = 56

$ws.Range("H25").Value = 150

$ws.Range("H31").Value = 3154.5334
$ws.Range("I31").Value = 2046.1333
$ws.Range("K31").Value = 2046.1333
$ws.Range("M31").Value = -1751.1333

$ws.Range("H34").Value = 3154.5334
$ws.Range("I34").Value = 2046.1333
$ws.Range("K34").Value = 2046.1333
$ws.Range("M34").Value = -1844.1333

$ws.Range("H86").Value = 15385.091
$ws.Range("I86").Value = 6347.8
$ws.Range("K86").Value = 6347.8
$ws.Range("M86").Value = -5224.8

$ws.Range("H89").Value = 15385.091
$ws.Range("I89").Value = 6347.8
$ws.Range("K89").Value = 31739
$ws.Range("M89").Value = -26123

$ws.Range("H132").Value = 4110.9033
$ws.Range("I132").Value = 2666.1177
$ws.Range("J132").Value = 5865.2856
$ws.Range("K132").Value = 7998.353099999999
$ws.Range("L132").Value = 17595.8568
$ws.Range("M132").Value = -5468.353099999999
$ws.Range("N132").Value = -22655.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 1000134.9
$ws.Range("I8").Value = 1000134.9
$ws.Range("K8").Value = 3000404.7
$ws.Range("M8").Value = -3000265.7

$ws.Range("H113").Value = 4115720.8
$ws.Range("J113").Value = 585.7143
$ws.Range("L113").Value = 1757.1429
$ws.Range("N113").Value = -6097.1429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()

$ws.Range("H97").Value = 21087.8
$ws.Range("I97").Value = 21087.8
$ws.Range("K97").Value = 21087.8
$ws.Range("M97").Value = -20591.8

$ws.Range("H122").Value = 6366.561
$ws.Range("I122").Value = 5864.643
$ws.Range("K122").Value = 17593.929
$ws.Range("M122").Value = -15143.929

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 334
$ws.Range("I9").Value = 223.33333
$ws.Range("K9").Value = 223.33333
$ws.Range("M9").Value = 0.6666700000000105

$ws.Range("H93").Value = 4749.75
$ws.Range("J93").Value = 2999.5
$ws.Range("L93").Value = 2999.5
$ws.Range("N93").Value = -5495.5

$ws.Range("H107").Value = 3498.25
$ws.Range("I107").Value = 3498.25
$ws.Range("K107").Value = 3498.25
$ws.Range("M107").Value = -1578.25

$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -55060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 11000
$ws.Range("J18").Value = 11000
$ws.Range("L18").Value = 11000
$ws.Range("N18").Value = -11346

$ws.Range("H54").Value = 37538.5
$ws.Range("J54").Value = 37538.5
$ws.Range("L54").Value = 37538.5
$ws.Range("N54").Value = -38578.5

$ws.Range("H57").Value = 98933.336
$ws.Range("J57").Value = 98933.336
$ws.Range("L57").Value = 98933.336
$ws.Range("N57").Value = -100441.336

$ws.Range("H62").Value = 112100
$ws.Range("I62").Value = 136375
$ws.Range("J62").Value = 15000
$ws.Range("K62").Value = 136375
$ws.Range("L62").Value = 15000
$ws.Range("M62").Value = -135751
$ws.Range("N62").Value = -16248

$ws.Range("H65").Value = 112100
$ws.Range("I65").Value = 136375
$ws.Range("J65").Value = 15000
$ws.Range("K65").Value = 681875
$ws.Range("L65").Value = 75000
$ws.Range("M65").Value = -678755
$ws.Range("N65").Value = -81240

$ws.Range("H100").Value = 700

$ws.Range("H107").Value = 102992.7
$ws.Range("I107").Value = 114333
$ws.Range("J107").Value = 930
$ws.Range("K107").Value = 342999
$ws.Range("L107").Value = 2790
$ws.Range("M107").Value = -341079
$ws.Range("N107").Value = -6630
